# ---------------------------------------------------------------------------
# AssureTestData.xlsx edit script
#
# Target shape (final tab order):
#   1. LoginTestSuccessFull     (unchanged content, selection -> A3)
#   2. LoginTestUnSuccessFull   (new header row, new column C text, extra row)
#   3. VerifySectionNames       (was "deals"; brand-new small lookup sheet)
#   4. VerifyButtonNames        (the actual former "deals" sheet, renamed &
#                                 refilled with the full button-name list)
#   5. SelectAbandonedVehicles  (brand-new sheet with a list validation)
#   6. tasks                    (unchanged content, simply moved to the end)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1) LoginTestSuccessFull -- cosmetic only: move the selection to A3
# ===========================================================================
$ws1 = $wb.Worksheets.Item("LoginTestSuccessFull")
$ws1.Activate()
$ws1.Range("A3").Select() | Out-Null

# ===========================================================================
# 2) LoginTestUnSuccessFull -- insert a title row, rename the error column,
#    swap in the new expected-error text, append one more data row
# ===========================================================================
$ws2 = $wb.Worksheets.Item("LoginTestUnSuccessFull")
$ws2.Activate()

$ws2.Rows.Item(1).Insert() | Out-Null
$ws2.Range("A1").Value = "LoginTestUnSuccessFull"

$ws2.Range("C2").Value = "ExpectedErrorMessage"

$errorMsg = "Please enter a valid email address/password. Please register before you login for the first time."
$ws2.Range("C3").Value = $errorMsg
$ws2.Range("C4").Value = $errorMsg

$ws2.Range("A5").Value = "abhilasha.jha@northgateps.com"
$ws2.Range("B5").Value = "N0rthg4t311"
$ws2.Range("C5").Value = $errorMsg

$ws2.Range("A2").Select() | Out-Null

# ===========================================================================
# 3) "deals" becomes "VerifyButtonNames" (keeps its historical sheetId) and
#    gets entirely new content -- the master list of button names.
# ===========================================================================
$deals = $wb.Worksheets.Item("deals")
$deals.Name = "VerifyButtonNames"
$deals.Activate()

# Row 2 (the styled header) is written before row 1 so the shared-string
# table fills in the same order the source workbook used.
$deals.Range("A2").Value = "ButtonNames"
$deals.Range("A2").Font.Bold = $true
$deals.Range("A2").Interior.Color = 65535

$deals.Range("A1").Value = "VerifyButtonNames"

$buttonNames = @(
    "Add OOA Address",
    "Configure SND",
    "Manage Contacts",
    "Manage Alerts",
    "Abandoned Vehicles",
    "Accidents",
    "Empty Homes",
    "Enquiries",
    "Environmental Permits",
    "Fixed Penalty Notice",
    "Food Registration",
    "HMO Details",
    "Infectious Disease",
    "Inspections Animal Feed",
    "Inspections Animal Health",
    "Inspections Environmental Permit",
    "Inspections Food Hygiene",
    "Inspections Food Scotland",
    "Inspections Food Standards",
    "Inspections Health And Safety",
    "Inspections HMO",
    "Inspections Licensing"
)
$row = 3
foreach ($name in $buttonNames) {
    $deals.Range("A" + $row).Value = $name
    $row = $row + 1
}

# Row 25 ("Inspections Non Routine") is filled in last, after the remaining
# rows 26-32, exactly mirroring the source edit order.
$buttonNamesTail = @(
    "Inspections Petroleum",
    "Inspections Primary Producers",
    "Inspections Trading Standards",
    "Licences",
    "Notices",
    "Prosecutions",
    "PSH Inspections"
)
$row = 26
foreach ($name in $buttonNamesTail) {
    $deals.Range("A" + $row).Value = $name
    $row = $row + 1
}
$deals.Range("A25").Value = "Inspections Non Routine"

$deals.Columns.Item(1).ColumnWidth = 27.6666666666667

# ===========================================================================
# 4) Brand-new "VerifySectionNames" sheet, inserted right before
#    "VerifyButtonNames" (takes over the freed-up sheetId).
# ===========================================================================
$vsn = $wb.Worksheets.Add($deals)
$vsn.Name = "VerifySectionNames"
$vsn.Activate()

$vsn.Range("A1").Value = "VerifySectionNames"
$vsn.Range("A2").Value = "SectionName"
$vsn.Range("A2").Font.Bold = $true
$vsn.Range("A2").Interior.Color = 65535
$vsn.Range("A3").Value = "General"
$vsn.Range("A4").Value = "Create Worksheets"

$vsn.Columns.Item(1).ColumnWidth = 16.5
$vsn.Columns.Item(2).ColumnWidth = 11

$vsn.Range("A8").Select() | Out-Null

# ===========================================================================
# 5) Brand-new "SelectAbandonedVehicles" sheet, inserted right after
#    "VerifyButtonNames", with a list data-validation on A3.
# ===========================================================================
$vbn = $wb.Worksheets.Item("VerifyButtonNames")
$sav = $wb.Worksheets.Add($null, $vbn)
$sav.Name = "SelectAbandonedVehicles"
$sav.Activate()

# Filled bottom-up, matching the source shared-string order (47, 48, 49).
$sav.Range("A3").Value = "ABVC"
$sav.Range("A2").Value = "AbandonedVehiclesType"
$sav.Range("A2").Font.Bold = $true
$sav.Range("A2").Interior.Color = 65535
$sav.Range("A1").Value = "selectAbandonedVehiclesType"

$sav.Range("A3").Validation.Add(3, 1, 1, '"ABVC,ABVD"') | Out-Null

$sav.Columns.Item(1).ColumnWidth = 25.3333333333333

$sav.Range("A3").Select() | Out-Null

# ===========================================================================
# 6) Move "tasks" to the end of the tab strip (content is unchanged).
# ===========================================================================
$tasks = $wb.Worksheets.Item("tasks")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tasks.Move($null, $lastSheet)

# ===========================================================================
# Final active tab: "SelectAbandonedVehicles" (5th tab, activeTab=4)
# ===========================================================================
$sav2 = $wb.Worksheets.Item("SelectAbandonedVehicles")
$sav2.Activate()
$sav2.Range("A3").Select() | Out-Null
